$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Move the existing "Group" output column (G20:G27) one column to the right (H20:H27),
#    preserving values, formulas and styles, to make room for the new "sum_weights" column.
$ws.Range("G20:G27").Copy($ws.Range("H20:H27"))

# 2) The new sum_weights column lands in what was blank space, so it should carry no
#    special formatting - clear whatever formatting the old "Group" column left behind.
$ws.Range("G20:G27").ClearFormats()

# 3) Header + formulas for the new column (each sums the weights feeding the weighted mean).
$ws.Range("G20").Value2 = "sum_weights"
$ws.Range("G21").Formula = "=F3+F16"
$ws.Range("G22").Formula = "=F4+F12+F17"
$ws.Range("G23").Formula = "=F5+F14"
$ws.Range("G24").Formula = "=F6+F8+F9"
$ws.Range("G25").Formula = "=F7+F10"
$ws.Range("G26").Formula = "=F11"
$ws.Range("G27").Formula = "=F13+F15"

# 4) The A:C labels for rows 21-27 switch from quoted-string values (shared with the
#    INPUT SAMPLE DATA table literal strings) to plain unquoted labels.
$ws.Range("A21").Value2 = "r000_100k"
$ws.Range("B21").Value2 = "r00_49"
$ws.Range("C21").Value2 = "white"

$ws.Range("A22").Value2 = "r000_100k"
$ws.Range("B22").Value2 = "r50plus"
$ws.Range("C22").Value2 = "black"

$ws.Range("A23").Value2 = "r000_100k"
$ws.Range("B23").Value2 = "r00_49"
$ws.Range("C23").Value2 = "black"

$ws.Range("A24").Value2 = "r100kplus"
$ws.Range("B24").Value2 = "r00_49"
$ws.Range("C24").Value2 = "aapi"

$ws.Range("A25").Value2 = "r100kplus"
$ws.Range("B25").Value2 = "r00_49"
$ws.Range("C25").Value2 = "aapi"

$ws.Range("A26").Value2 = "negative"
$ws.Range("B26").Value2 = "r50plus"
$ws.Range("C26").Value2 = "hispanic"

$ws.Range("A27").Value2 = "r000_100k"
$ws.Range("B27").Value2 = "r50plus"
$ws.Range("C27").Value2 = "aian"

# 5) These cells carried a redundant "no border, no fill" style left over from earlier
#    edits; clear it so they fall back to the sheet's plain default formatting.
$ws.Range("E4").ClearFormats()
$ws.Range("E7").ClearFormats()
$ws.Range("E8").ClearFormats()
$ws.Range("E9").ClearFormats()
$ws.Range("E27").ClearFormats()
